# Insert a new data row before the current row 682 (Camote / "1a (cosecha)" /
# Region de O'Higgins, dated 2023-02-27 = serial 44984). Excel's Insert shifts
# every subsequent row down by one, which reproduces the rest of the diff
# (all rows 682-773 become 683-774, with the former last row, 773, becoming
# the new last row, 774) without any further edits needed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(682).Insert()

$ws.Range("A682").Value = 10
$ws.Range("B682").Value = "Vega Modelo de Temuco"
$ws.Range("C682").Value = "La Araucanía"
$ws.Range("D682").Value = 44984
$ws.Range("E682").Value = 9
$ws.Range("F682").Value = 100112045
$ws.Range("G682").Value = "Zapallo"
$ws.Range("H682").Value = "Camote"
$ws.Range("I682").Value = "1a (cosecha)"
$ws.Range("J682").Value = 1800
$ws.Range("K682").Value = 500
$ws.Range("L682").Value = 700
$ws.Range("M682").Value = 589
$ws.Range("N682").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O682").Value = "Región de O'Higgins"
$ws.Range("P682").Value = 589
$ws.Range("Q682").Value = 1
$ws.Range("R682").Value = "Hortaliza"
